$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.179.89'
$ws.Range('E2').Value = '  +0.01%  '

$ws.Range('D3').Value = '3.772.99'
$ws.Range('E3').Value = '  -1.45%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.89%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.10%  '

$ws.Range('D7').Value = '3.768.87'
$ws.Range('E7').Value = '  -1.39%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.523'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.62%  '

$ws.Range('E10').Value = '  -2.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.48'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.79%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.451'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E13').Value = '  -2.28%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.76%  '

$ws.Range('D15').Value = '4.399.80'
$ws.Range('E15').Value = '  -1.40%  '

$ws.Range('D16').Value = '3.766.35'
$ws.Range('E16').Value = '  -1.31%  '

$ws.Range('D17').Value = '68.172.03'
$ws.Range('E17').Value = '  +0.14%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.15%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.08%  '

$ws.Range('E20').Value = '  -0.36%  '

$ws.Range('E21').Value = '  +1.00%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '469.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.705'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.27%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.27%  '

$ws.Range('E25').Value = '  -6.52%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.01%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.38%  '

$ws.Range('E28').Value = '  -2.00%  '

$ws.Range('E29').Value = '  +0.06%  '

$ws.Range('D30').Value = '3.913.50'
$ws.Range('E30').Value = '  -1.34%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.05%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.45'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.77%  '

$ws.Range('E33').Value = '  -1.79%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.28%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.22%  '

$ws.Range('D37').Value = '3.723.22'

$ws.Range('E38').Value = '  -3.44%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.49'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.39%  '

$ws.Range('E40').Value = '  -0.69%  '

$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.59%  '

$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.996'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.90%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.22%  '

$ws.Range('E44').Value = '  -3.32%  '

$ws.Range('B46').Value = 'Cosmos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.64%  '

$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.13%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '400.85'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.40%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.26%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.19%  '
